$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Toggle switches grid (A1:C5) -----------------------------------------
# Before: B1=1, C2=1, C4=1 (others blank)
# After:  A2=1, B3=1, A4=1, C4=1 (B1 and C2 cleared)
$ws.Range("B1").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("A2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("A4").Value = 1

# --- New lookup-table rows 30:32 (I/J/K columns) ---------------------------
# Extend the hex-digit -> display-character lookup table used by the game.
$ws.Range("I30").Value = "b"
$ws.Range("J30").Value = 182

$ws.Range("I31").Value = "'1."
$ws.Range("J31").Formula = "=5+8"

$ws.Range("I32").Value = "h"
$ws.Range("J32").Value = 150

# Re-apply the DEC2HEX formula across the whole K column so it becomes one
# shared formula group spanning K2:K32 (matches rows 2-29 plus the new 30-32).
$ws.Range("K2:K32").Formula = "=DEC2HEX(J2)"

# --- New "Q"/"R" helper columns --------------------------------------------
# Q holds a display-character (mirrors the I column); R looks up its hex
# code from the K column - used to build the start-up / score strings.
$ws.Range("Q5").Value = "b"
$ws.Range("R5").Formula = "=K30"

$ws.Range("Q6").Value = "I"
$ws.Range("R6").Formula = "=K28"

$ws.Range("Q7").Value = "n"
$ws.Range("R7").Formula = "=K20"

$ws.Range("Q8").Value = "A"
$ws.Range("R8").Formula = "=K12"

$ws.Range("Q9").Value = "r"
$ws.Range("Q9").Font.Bold = $true
$ws.Range("R9").Formula = "=K26"

$ws.Range("Q10").Value = "Y"
$ws.Range("R10").Formula = "=K19"

$ws.Range("Q12").Value = "g"
$ws.Range("R12").Formula = "=K23"

$ws.Range("Q13").Value = "A"
$ws.Range("R13").Formula = "=K12"

$ws.Range("Q14").Value = "m"
$ws.Range("R14").Formula = "=K22"

$ws.Range("Q15").Value = "E"
$ws.Range("R15").Formula = "=K16"

$ws.Range("Q17").Value = "'1."
$ws.Range("R17").Formula = "=K31"

$ws.Range("Q18").Value = 0
$ws.Range("R18").Formula = "=K2"

$ws.Range("Q24").Value = "O"
$ws.Range("R24").Formula = "=K24"

$ws.Range("Q25").Value = "F"
$ws.Range("R25").Formula = "=K17"

$ws.Range("Q26").Value = "F"
$ws.Range("R26").Formula = "=K17"

$ws.Range("Q28").Value = "O"
$ws.Range("R28").Formula = "=K24"

$ws.Range("Q29").Value = "n"
$ws.Range("R29").Formula = "=K20"

# Remaining new R cells (no matching Q cell) are plain zero placeholders.
$ws.Range("R1").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("R16").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("R21").Value = 0
$ws.Range("R22").Value = 0

# --- Selection --------------------------------------------------------------
$ws.Range("R24:R26").Select() | Out-Null
